# Added results for partial model
# Populate the "IPC PO" (predicted) column together with the derived
# DELTA / DELTA^2 columns and the summary TOTAL / MSE rows for the
# partial-model sliding-window results (window 7).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New predicted values (column C) for data rows 2-51
$cVals = @("29.64127922058105", "29.90331840515137", "29.90265464782715", "29.92182731628418", "29.99698448181152", "30.15762138366699", "30.14980888366699", "30.26547050476074", "30.23031425476074", "30.44497489929199", "30.40680503845215", "30.44501304626465", "30.42711448669434", "30.66121482849121", "30.80583000183105", "31.21635246276855", "31.35795402526855", "31.45592308044434", "31.29449272155762", "31.48099327087402", "32.05437088012695", "32.60719680786133", "32.54343032836914", "32.63210678100586", "32.83679580688477", "33.05073165893555", "32.78688430786133", "33.60528182983398", "33.57837295532227", "33.80495834350586", "34.3074836730957", "35.01227951049805", "35.66598129272461", "35.98637771606445", "35.96036148071289", "36.60099411010742", "37.27934646606445", "38.10490798950195", "38.57924270629883", "39.20572280883789", "39.69402694702148", "39.77173233032227", "39.85319900512695", "40.23218154907227", "40.4200325012207", "40.79880142211914", "40.99710464477539", "41.71234130859375", "41.72422790527344", "43.96305465698242")

# Corresponding DELTA values (column D = C - B) for data rows 2-51
$dVals = @("-0.1087207794189453", "0.06331840515136378", "0.09265464782714616", "0.001827316284177982", "0.01698448181151946", "0.1176213836670001", "-0.06019111633301577", "0.04547050476074332", "-0.1496857452392533", "0.004974899291994461", "-0.07319496154785554", "-0.2449869537353493", "-0.3228855133056641", "-0.2787851715087868", "-0.1441699981689482", "0.1963524627685587", "0.2379540252685501", "0.1759230804443348", "-0.08550727844237827", "-0.09900672912597486", "0.4043708801269474", "0.7271968078613327", "0.2634303283691395", "0.1821067810058565", "-0.01320419311522869", "0.1507316589355412", "-0.3131156921386662", "0.2052818298339787", "-0.1216270446777372", "-0.2950416564941349", "-0.09251632690430256", "0.1122795104980412", "0.3659812927246122", "0.2863777160644503", "-0.3396385192871065", "-0.1990058898925753", "-0.02065353393554403", "0.2049079895019474", "0.07924270629882812", "0.3057228088378849", "0.2940269470214787", "-0.1282676696777401", "-0.2468009948730412", "-0.3678184509277287", "-0.4799674987793026", "-0.4011985778808693", "-0.5028953552246094", "-0.08765869140624716", "-0.4757720947265653", "1.263054656982419")

# Corresponding DELTA^2 values (column E = D^2) for data rows 2-51
$eVals = @("0.01182020787746296", "0.004009220430912251", "0.008584883763972481", "3.339084802422028e-06", "0.0002884726224058353", "0.01383478989573965", "0.003622970485414638", "0.002067566803196782", "0.02240582232783063", "2.474962296548699e-05", "0.005357502395992051", "0.06001860750052617", "0.1042550547026622", "0.07772117185318367", "0.02078498837203452", "0.03855428963527822", "0.0566221181415058", "0.03094893023302389", "0.007311494666622407", "0.009802332412224158", "0.1635158086946421", "0.528815197363712", "0.06939553790467266", "0.03316287968831499", "0.0001743507158242527", "0.02272003300546031", "0.09804143666347598", "0.04214062965998659", "0.01479313799704029", "0.08704957906680312", "0.008559270743863777", "0.01260668847767974", "0.1339423066243783", "0.08201219625829091", "0.1153543237835382", "0.03960334421193579", "0.000426568464026669", "0.04198728416173021", "0.006279406501562335", "0.09346643584372594", "0.08645184557477144", "0.01645259508455784", "0.0609107310703229", "0.135290412842874", "0.2303687998844598", "0.160960298893632", "0.252903738306486", "0.00768404617905567", "0.2263590861205039", "1.595307066524976")

$firstRow = 2
for ($i = 0; $i -lt $cVals.Length; $i++) {
    $r = $firstRow + $i
    $ws.Cells.Item($r, 3).Value = [double]$cVals[$i]
    $ws.Cells.Item($r, 4).Value = [double]$dVals[$i]
    $ws.Cells.Item($r, 5).Value = [double]$eVals[$i]
}

# A handful of RO (column B) values were re-saved with a tiny floating
# point representation change (double round-trip) as part of this export.
$bRows = @(7, 8, 22, 26, 27, 28, 29, 31, 32, 33, 39, 41, 42, 43, 44, 45, 46, 47)
$bVals = @("30.03999999999999", "30.21000000000001", "31.65000000000001", "32.84999999999999", "32.90000000000001", "33.09999999999999", "33.40000000000001", "34.09999999999999", "34.40000000000001", "34.90000000000001", "37.90000000000001", "38.90000000000001", "39.40000000000001", "39.90000000000001", "40.09999999999999", "40.59999999999999", "40.90000000000001", "41.20000000000001")

for ($i = 0; $i -lt $bRows.Length; $i++) {
    $ws.Cells.Item($bRows[$i], 2).Value = [double]$bVals[$i]
}

# Summary rows: TOTAL (52) and MSE (53)
$ws.Range("C52").Value = [double]"0.1454766845702764"
$ws.Range("E52").Value = [double]"4.84477354914006"
$ws.Range("E53").Value = [double]"0.0968954709828012"
